# Applies two edits described by the target diff:
#  1) "Hugo Pelayo Aseko" (with spell-check proofErr wrapper around "Aseko")
#     becomes "Hugo Pelayo" - i.e. the second run/word + proofErr markers are removed.
#  2) The four runs forming the "Es la situacion ... horario de" paragraph are
#     merged into a single run, and the paragraph's <w:pPr><w:spacing w:after="0"/></w:pPr>
#     is removed.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wNs + '><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Change 1: "Hugo Pelayo Aseko" -> "Hugo Pelayo" ------------------------

$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Hugo Pelayo Aseko", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $target1 = $d.Range($rng1.Start, $rng1.End)
    $innerXml1 = '<w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Hugo Pelayo</w:t></w:r></w:p>'
    $null = $target1.InsertXML($pkgHeader + $innerXml1 + $pkgFooter)
} else {
    Write-Host "WARNING: 'Hugo Pelayo Aseko' not found"
}

# --- Change 2: merge the desconexion digital paragraph's runs --------------

$mergedText = "Es la situaci" + [char]0xF3 + "n en la cual el trabajador no est" + [char]0xE1 + " en el deber de responder al tel" + [char]0xE9 + "fono, a los correos electr" + [char]0xF3 + "nicos o mensajes profesionales de cualquier otro tipo, etc. fuera de su horario de"

$rng2 = $d.Content
$found2 = $rng2.Find.Execute($mergedText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $target2 = $d.Range($rng2.Start, $rng2.End)
    $innerXml2 = '<w:p><w:r><w:t>' + $mergedText + '</w:t></w:r></w:p>'
    $null = $target2.InsertXML($pkgHeader + $innerXml2 + $pkgFooter)
} else {
    Write-Host "WARNING: desconexion paragraph not found"
}
